$wb = $excel.ActiveWorkbook
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "DUA"
$newSheet.Move($wb.Worksheets.Item("authors"))
$dua = $wb.Worksheets.Item("DUA")
$dua.Activate()

$dua.Range("D2").Value = "x"
$dua.Range("B2").Validation.Add(3, 1, 1, "=`$D`$2:`$D`$6")
$dua.Range("B2").Validation.InputMessage = "select from drop down menu"
$dua.Range("B2").Validation.ShowInput = $true
Write-Output "done"
